$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update unit price column (G) for a few rows, per diff:
# G2: 3.15 -> 2.95
# G4: 3.25 -> 2.95
# G5: 3.30 -> 2.95
# G6: 3.80 -> 3.15
$ws.Range("G2").Value = 2.95
$ws.Range("G4").Value = 2.95
$ws.Range("G5").Value = 2.95
$ws.Range("G6").Value = 3.15
